$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone the S:T block formatting onto the new X:Y block (row 3 has no T3, so copy it alone) ---
$ws.Range("S3").Copy()
$ws.Range("X3").PasteSpecial(-4122)
$ws.Range("S4:T57").Copy()
$ws.Range("X4").PasteSpecial(-4122)

# --- Header row 3: swap the "3Y" title from S3 to X3, put new "2Y" title into S3 ---
$ws.Range("S3").Value = "RF-100 (superdataset-24-f 2Y.csv + extrapol)"
$ws.Range("X3").Value = "RF-100 (superdataset-24-f 3Y.csv + extrapol)"

# --- Header row 4: new block gets the same "test (R2)" label ---
$ws.Range("Y4").Value = "test (R2)"

# --- New index column X (1..50), mirroring C/H/N/S columns ---
$ws.Range("X5").Value = 1
$ws.Range("X6").Formula = "=X5+1"
$ws.Range("X7:X54").Formula = "=X6+1"

# --- S/T block gets freshly computed "2Y" test-R2 values ---
$ws.Range("T5").Value = 0.7475993806200267
$ws.Range("T6").Value = 0.7527221679329714
$ws.Range("T7").Value = 0.8129006274528505
$ws.Range("T8").Value = 0.763596592812755
$ws.Range("T9").Value = 0.7545887963884832
$ws.Range("T10").Value = 0.8038516673740239
$ws.Range("T11").Value = 0.8153705922463242
$ws.Range("T12").Value = 0.8042547093447852
$ws.Range("T13").Value = 0.7576941490745841
$ws.Range("T14").Value = 0.7466729504081231
$ws.Range("T15").Value = 0.7847900633381837
$ws.Range("T16").Value = 0.774985082871346
$ws.Range("T17").Value = 0.7668050636127745
$ws.Range("T18").Value = 0.7954786090034827
$ws.Range("T19").Value = 0.8025392379716025
$ws.Range("T20").Value = 0.7581103740855037
$ws.Range("T21").Value = 0.8095949650875207
$ws.Range("T22").Value = 0.7824322314568165
$ws.Range("T23").Value = 0.790655312957717
$ws.Range("T24").Value = 0.7511117733443164
$ws.Range("T25").Value = 0.7933349338637159
$ws.Range("T26").Value = 0.7431987341233699
$ws.Range("T27").Value = 0.8048189704153046
$ws.Range("T28").Value = 0.7783134716136362
$ws.Range("T29").Value = 0.7325044742101615
$ws.Range("T30").Value = 0.7992896977623103
$ws.Range("T31").Value = 0.8014504030103513
$ws.Range("T32").Value = 0.7920319778958196
$ws.Range("T33").Value = 0.7454934779298727
$ws.Range("T34").Value = 0.7835693262431689
$ws.Range("T35").Value = 0.7816576623573809
$ws.Range("T36").Value = 0.8129048890288546
$ws.Range("T37").Value = 0.8036130733916196
$ws.Range("T38").Value = 0.7785015148210379
$ws.Range("T39").Value = 0.7635795733540454
$ws.Range("T40").Value = 0.8271959110509826
$ws.Range("T41").Value = 0.7876871863039125
$ws.Range("T42").Value = 0.8072189130432519
$ws.Range("T43").Value = 0.8068480169465223
$ws.Range("T44").Value = 0.7786592700395407
$ws.Range("T45").Value = 0.7844005819567024
$ws.Range("T46").Value = 0.7751775756956916
$ws.Range("T47").Value = 0.731721785307234
$ws.Range("T48").Value = 0.7873971117873677
$ws.Range("T49").Value = 0.7698124543944495
$ws.Range("T50").Value = 0.7431058421406191
$ws.Range("T51").Value = 0.8032203794821117
$ws.Range("T52").Value = 0.7819228097221358
$ws.Range("T53").Value = 0.7402717713161295
$ws.Range("T54").Value = 0.770664507629979

# --- New X/Y block gets the values that used to live in S/T ("3Y" results) ---
$ws.Range("Y5").Value = 0.7180455882572573
$ws.Range("Y6").Value = 0.7218538185963723
$ws.Range("Y7").Value = 0.6981084741413246
$ws.Range("Y8").Value = 0.7182924440651435
$ws.Range("Y9").Value = 0.6504707959236231
$ws.Range("Y10").Value = 0.635028559524275
$ws.Range("Y11").Value = 0.713213169608098
$ws.Range("Y12").Value = 0.6901136174071281
$ws.Range("Y13").Value = 0.7056104872858253
$ws.Range("Y14").Value = 0.6681121073668193
$ws.Range("Y15").Value = 0.6855759496883969
$ws.Range("Y16").Value = 0.6851485124334683
$ws.Range("Y17").Value = 0.62803686401784
$ws.Range("Y18").Value = 0.6947538021994439
$ws.Range("Y19").Value = 0.6627521297507943
$ws.Range("Y20").Value = 0.6684050381874103
$ws.Range("Y21").Value = 0.6433350510301314
$ws.Range("Y22").Value = 0.6954556236691127
$ws.Range("Y23").Value = 0.6844023612177231
$ws.Range("Y24").Value = 0.6639150797752684
$ws.Range("Y25").Value = 0.66425350778334
$ws.Range("Y26").Value = 0.718111569087034
$ws.Range("Y27").Value = 0.58926158900662
$ws.Range("Y28").Value = 0.6880634305504858
$ws.Range("Y29").Value = 0.6624016276088496
$ws.Range("Y30").Value = 0.6296252087808936
$ws.Range("Y31").Value = 0.7042949611769824
$ws.Range("Y32").Value = 0.7346055817943136
$ws.Range("Y33").Value = 0.6695388415189456
$ws.Range("Y34").Value = 0.6802785261892927
$ws.Range("Y35").Value = 0.6993139246063016
$ws.Range("Y36").Value = 0.7231808762501936
$ws.Range("Y37").Value = 0.6764350080179606
$ws.Range("Y38").Value = 0.6368008590747593
$ws.Range("Y39").Value = 0.7574871794856763
$ws.Range("Y40").Value = 0.6527747799660627
$ws.Range("Y41").Value = 0.6705588084872671
$ws.Range("Y42").Value = 0.6265846458336437
$ws.Range("Y43").Value = 0.7058349624964528
$ws.Range("Y44").Value = 0.65343318008461
$ws.Range("Y45").Value = 0.6372088498636099
$ws.Range("Y46").Value = 0.7059213929907977
$ws.Range("Y47").Value = 0.6260447049353088
$ws.Range("Y48").Value = 0.6662733373015717
$ws.Range("Y49").Value = 0.6629240671107928
$ws.Range("Y50").Value = 0.6507557408950848
$ws.Range("Y51").Value = 0.707939323290021
$ws.Range("Y52").Value = 0.6864707257005609
$ws.Range("Y53").Value = 0.6693298947890928
$ws.Range("Y54").Value = 0.7274695469532186

# --- avg / SD rows ---
$ws.Range("X56").Value = "avg"
$ws.Range("Y56").Formula = "=AVERAGE(Y5:Y54)"
$ws.Range("X57").Value = "SD"
$ws.Range("Y57").Formula = "=_xlfn.STDEV.S(Y5:Y54)"

# --- Update selection / scroll position to match the saved view ---
$ws.Range("U71").Select()
